$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, shifting existing rows 21-31 down to 22-32
$ws.Rows(21).Insert()

# Populate the newly inserted row 21 with the new weekly price record
$ws.Cells.Item(21, 1).Value = 8
$ws.Cells.Item(21, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(21, 3).Value = "Coquimbo"
$ws.Cells.Item(21, 4).Value = 44596
$ws.Cells.Item(21, 5).Value = 4
$ws.Cells.Item(21, 6).Value = 100114007
$ws.Cells.Item(21, 7).Value = "Jengibre"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 500
$ws.Cells.Item(21, 11).Value = 16000
$ws.Cells.Item(21, 12).Value = 17000
$ws.Cells.Item(21, 13).Value = 16500
$ws.Cells.Item(21, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(21, 15).Value = "Perú"
$ws.Cells.Item(21, 16).Value = 1269
$ws.Cells.Item(21, 17).Value = 13
$ws.Cells.Item(21, 18).Value = "Hortaliza"
